$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Occurrences": clean up duplicate/erroneous occurrence rows
# ---------------------------------------------------------------------
$occ = $wb.Worksheets.Item("Occurrences")

# Row 2 (UNCEN-2000HP-HS002-CA001-VE001 / Phalanger orientalis):
# the individual count was wrong and the occurrenceStatus cell had
# drifted onto the wrong shared string - fix both.
$occ.Range("M2").Value = 3
$occ.Range("P2").Value = "1999-09-16/1999-10-07"

# Row 3 used to hold a duplicate "CA001-VE002" Phalanger orientalis
# record - that whole record is removed and replaced with the data
# that used to live in row 5 (PW001 / Spilocuscus maculatus), which
# collapses the trailing duplicate rows.
$occ.Range("A3").Value = "UNCEN-2000HP-HS002-PW001"
$occ.Range("B3").Value = "UNCEN-2000HP-HS002-PW001-VE001"
$occ.Range("F3").Value = "Spilocuscus maculatus"
$occ.Range("A3").Copy()
$occ.Range("G3").PasteSpecial(-4122)   # xlPasteFormats - drop G3's old style, matching G5's
$occ.Range("P3").Value = "1999-09-16/1999-10-07"

# Row 4 (old CA001-VE003 duplicate) is blanked out entirely, leaving
# only the styled-but-empty placeholder cells used by the rest of the
# sheet's template rows.
$occ.Range("A4:E4").ClearContents()
$occ.Range("F4").ClearContents()
$occ.Range("G4").Clear()
$occ.Range("I4:J4").ClearContents()
$occ.Range("K4:M4").ClearContents()
$occ.Range("P4").ClearContents()
$occ.Range("I6").Copy()
$occ.Range("I4").PasteSpecial(-4122)
$occ.Range("J6").Copy()
$occ.Range("J4").PasteSpecial(-4122)

# Row 5 (old PW001-VE001 record, now merged into row 3 above) is
# blanked out the same way.
$occ.Range("A5:E5").ClearContents()
$occ.Range("F5").ClearContents()
$occ.Range("G5").Clear()
$occ.Range("I5:J5").ClearContents()
$occ.Range("K5:M5").ClearContents()
$occ.Range("P5").ClearContents()
$occ.Range("I6").Copy()
$occ.Range("I5").PasteSpecial(-4122)
$occ.Range("J6").Copy()
$occ.Range("J5").PasteSpecial(-4122)

# View/column tidy-up that came along with the data cleanup.
$occ.Range("D10").Select()
$occ.Columns.Item(1).ColumnWidth = 33

# ---------------------------------------------------------------------
# Sheet "Sampling Events": widen the eventID column to fit the
# shortened set of sampling-event identifiers.
# ---------------------------------------------------------------------
$se = $wb.Worksheets.Item("Sampling Events")
$se.Columns.Item(2).ColumnWidth = 29
